$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values
$ws.Range("B2").Value = 183
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 143

# Add new rows of data, copying style from A3 (which already has the needed style)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 65

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 44

$ws.Range("A6").Value = 1
$ws.Range("B6").Value = 34

# Apply the same style (used by A2/A3) to the new A4:A6 cells
$ws.Range("A2:A3").Copy() | Out-Null
$ws.Range("A4:A6").PasteSpecial(-4122) | Out-Null
